# Update countries & provincias Spain
#
# The source feed re-ranked a few countries (Armenia and Georgia moved up
# in the ranking) which cascades into the rows below them, and refreshed
# a couple of independent country totals (Noruega, India) plus the
# "last updated" timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." banner: 07:16 -> 07:46
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 07:46"

# Noruega (row 17): refreshed totals
$ws.Range("B17").Value = 1967
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 1959

# India (row 49): refreshed totals
$ws.Range("B49").Value = 275
$ws.Range("C49").Value = 26
$ws.Range("E49").Value = 247

# Armenia jumps above Kuwait/Colombia/Argentina/San Marino/
# Emiratos Arabes Unidos/Eslovaquia, pushing each of them down one row.

# Row 59: Kuwait -> Armenia (brand-new totals)
$ws.Range("A59").Value = "Armenia"
$ws.Range("B59").Value = 160
$ws.Range("C59").Value = 24
$ws.Range("D59").Value = 1
$ws.Range("E59").Value = 159
$ws.Range("F59").Value = 2

# Row 60: Colombia -> Kuwait (carries Kuwait's previous totals)
$ws.Range("A60").Value = "Kuwait"
$ws.Range("B60").Value = 159
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 22
$ws.Range("E60").Value = 137
$ws.Range("F60").Value = 5

# Row 61: Argentina -> Colombia (carries Colombia's previous totals)
$ws.Range("A61").Value = "Colombia"
$ws.Range("C61").Value = 13
$ws.Range("D61").Value = 1
$ws.Range("E61").Value = 157
$ws.Range("H61").Value = 0

# Row 62: San Marino -> Argentina (carries Argentina's previous totals)
$ws.Range("A62").Value = "Argentina"
$ws.Range("B62").Value = 158
$ws.Range("D62").Value = 3
$ws.Range("E62").Value = 152
$ws.Range("F62").Value = 0
$ws.Range("H62").Value = 3

# Row 63: Emiratos Arabes Unidos -> San Marino (carries its previous totals)
$ws.Range("A63").Value = "San Marino"
$ws.Range("B63").Value = 151
$ws.Range("D63").Value = 4
$ws.Range("E63").Value = 133
$ws.Range("F63").Value = 12
$ws.Range("H63").Value = 14

# Row 64: Eslovaquia -> Emiratos Arabes Unidos (carries its previous totals)
$ws.Range("A64").Value = "Emiratos Arabes Unidos"
$ws.Range("B64").Value = 140
$ws.Range("D64").Value = 31
$ws.Range("E64").Value = 107
$ws.Range("H64").Value = 2

# Row 65: Armenia (old) -> Eslovaquia (carries its previous totals)
$ws.Range("A65").Value = "Eslovaquia"
$ws.Range("B65").Value = 137
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 137

# Georgia jumps above Senegal, pushing it down one row.

# Row 99: Senegal -> Georgia (brand-new totals)
$ws.Range("A99").Value = "Georgia"
$ws.Range("C99").Value = 3
$ws.Range("D99").Value = 1
$ws.Range("E99").Value = 46
$ws.Range("F99").Value = 1

# Row 100: Georgia (old) -> Senegal (carries Senegal's previous totals)
$ws.Range("A100").Value = "Senegal"
$ws.Range("B100").Value = 47
$ws.Range("D100").Value = 5
$ws.Range("E100").Value = 42
$ws.Range("F100").Value = 0
